$d = $word.ActiveDocument

# --- Edit 1: non_payment_hardship / non_payment_action -> nonpayment_appl_expired ---
# Replaces "non_payment_hardship and non_payment_action" (spanning 3 runs) with a single
# "nonpayment_appl_expired" so the checkbox reads output_checkbox(nonpayment_appl_expired)
$found1 = $d.Content.Find.Execute(
    "non_payment_hardship and non_payment_action",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "nonpayment_appl_expired", 2)

# --- Edit 2: eviction_reason_nofault -> nofault_judgment_period_exp ---
$found2 = $d.Content.Find.Execute(
    "output_checkbox(eviction_reason_nofault",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "output_checkbox(nofault_judgment_period_exp", 2)

# --- Edit 3: eviction_reason_139 -> eviction_reason_139 and judgment_outcome == "landlord" ---
# Uses InsertAfter (rather than Find/Replace's replacement text) so the straight quotes
# around "landlord" are not mangled into curly quotes by AutoFormat/AutoCorrect.
$r3 = $d.Content
$found3 = $r3.Find.Execute(
    "output_checkbox(eviction_reason_139",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$ins3 = $d.Range($r3.End, $r3.End)
$ins3.InsertAfter(" and judgment_outcome == " + [char]34 + "landlord" + [char]34)
